# Update "想去人数" (interest/attendance count) figures in the F column
# across the "展览", "本地生活" and "全部类型" sheets, matching the
# upstream data refresh recorded in the commit "Update gh-pages to output
# generated at 7921097".

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" (Exhibitions) ---
$ws = $wb.Worksheets.Item("展览")
$ws.Range("F7").Value = 9447
$ws.Range("F10").Value = 685
$ws.Range("F11").Value = 1920
$ws.Range("F12").Value = 43
$ws.Range("F13").Value = 480
$ws.Range("F14").Value = 2572
$ws.Range("F16").Value = 3870
$ws.Range("F18").Value = 135
$ws.Range("F20").Value = 205
$ws.Range("F22").Value = 17
$ws.Range("F26").Value = 550
$ws.Range("F27").Value = 2130
$ws.Range("F28").Value = 1090
$ws.Range("F30").Value = 463
$ws.Range("F33").Value = 139
$ws.Range("F35").Value = 127

# --- Sheet "本地生活" (Local life) ---
$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F2").Value = 185
$ws.Range("F3").Value = 966

# --- Sheet "全部类型" (All types) ---
$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F2").Value = 185
$ws.Range("F4").Value = 966
$ws.Range("F10").Value = 9447
$ws.Range("F13").Value = 685
$ws.Range("F14").Value = 1920
$ws.Range("F15").Value = 43
$ws.Range("F16").Value = 480
$ws.Range("F18").Value = 2572
$ws.Range("F20").Value = 3870
$ws.Range("F22").Value = 135
$ws.Range("F24").Value = 205
$ws.Range("F26").Value = 17
$ws.Range("F31").Value = 550
$ws.Range("F32").Value = 2130
$ws.Range("F33").Value = 1090
$ws.Range("F35").Value = 463
$ws.Range("F38").Value = 139
$ws.Range("F40").Value = 127
